$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of column letter -> array of new values for rows 2..25
$colData = @{}
$colData["B"] = @(2.07661134372313, 1.940982682918502, 1.858873804843142, 1.825705081035835, 1.820214967828861, 1.858425303360548, 2.029602781777442, 2.374671560740751, 2.634122376775622, 2.753487500258473, 2.798883907415927, 2.78909824212127, 2.757218357154954, 2.737716541228394, 2.626348778790145, 2.558373337961996, 2.519401584964612, 2.506227962688911, 2.565596387219045, 2.766576925782033, 2.899069761881208, 2.828250551649376, 2.562330509736455, 2.280295163257222)
$colData["C"] = @(0.4821953700043196, 0.443013500792091, 0.4192009147496378, 0.4095580047718386, 0.4079604641194408, 0.4190706214428417, 0.468634090284354, 0.567813695820746, 0.6419552937402386, 0.6759754768783637, 0.6889011829067044, 0.6861154738400046, 0.6770380154686677, 0.6714834371193774, 0.6397379627795203, 0.6203387663260287, 0.6092083661794163, 0.6054445009119149, 0.6224009908617631, 0.6797031107390126, 0.717404433890124, 0.6972592525116283, 0.6214685892697958, 0.540763877812708)
$colData["D"] = @(0.6647245498891436, 0.6552736699228205, 0.6498472391969017, 0.6477302368544144, 0.6473843964304535, 0.6498183071493884, 0.6613875381064815, 0.6870797433008988, 0.7078160465485439, 0.7176599281451672, 0.7214471096501711, 0.7206288181841103, 0.7179703064730631, 0.7163496554316851, 0.7071810300924994, 0.7016619090582026, 0.6985260991960729, 0.6974709915405981, 0.7022454264047155, 0.7187495569641555, 0.7298830844020472, 0.7239090010752705, 0.7019815023554372, 0.6798047277929129)
$colData["E"] = @(0.2709984499433773, 0.2662864753111407, 0.2635539597482719, 0.2624806502429493, 0.262304851424112, 0.2635393221418667, 0.2693403021818312, 0.2820001651109081, 0.2920988158631204, 0.296869423546525, 0.2987015841152072, 0.2983058517377515, 0.2970196414079354, 0.2962351457467136, 0.2917906201110299, 0.289109460668989, 0.2875839447192448, 0.2870702781779642, 0.2893931533330374, 0.2973967349864637, 0.3027770684382389, 0.2998917242885, 0.2892648463147296, 0.2784363683575606)
$colData["G"] = @(2.045243200238332, 2.00779472421118, 1.986296759856145, 1.977909403906409, 1.976539150441766, 1.986182137555744, 2.032018835524809, 2.133906850354037, 2.216281013717577, 2.255431906987866, 2.270502170192515, 2.267245584395823, 2.256666823938929, 2.250218989954845, 2.21375646746182, 2.191820124899891, 2.179360705831442, 2.175169160398781, 2.194138932081188, 2.259767395436569, 2.304087084250767, 2.280301091656838, 2.193090125919781, 2.105037997985534)
$colData["H"] = @(1.658982787693475, 1.649487226274772, 1.644661129361339, 1.642945834149856, 1.642676160307559, 1.644636980000087, 1.655499595076293, 1.684824399769468, 1.71134258350159, 1.724504665527917, 1.729648255001308, 1.728533379520542, 1.724924629201098, 1.722734964855391, 1.71050463946969, 1.703284143988384, 1.699234448445367, 1.697880998500267, 1.704042072560611, 1.725980269554924, 1.741247884878504, 1.733013736041926, 1.703699097121159, 1.67602392367013)
$colData["J"] = @(0.1409426208235516, 0.1377665589584325, 0.1359051221321721, 0.1351687364383736, 0.135047794994783, 0.1358951014058505, 0.1398290270486058, 0.14825375156164, 0.1548866232292738, 0.158002632135549, 0.1591969486833165, 0.158939090351339, 0.1581006007134533, 0.1575888754451853, 0.1546849841620457, 0.152928929449871, 0.1519281718263557, 0.1515909212434394, 0.15311490301508, 0.1583464946451727, 0.1618493716042622, 0.1599721048381042, 0.1530307969011204, 0.1458975536202871)
$colData["N"] = @(1.540212862536709, 1.560852772425484, 1.574170500552352, 1.579759233052759, 1.58069698732686, 1.574245218288567, 1.547195402628994, 1.499287098368466, 1.467247385873932, 1.453363532262934, 1.448206023694901, 1.449312327734575, 1.452937213102707, 1.455170602196247, 1.468168688874457, 1.476320054164916, 1.48107352276795, 1.48269411403335, 1.475445593787327, 1.451869777313377, 1.437044821080242, 1.444903609557912, 1.475840728244194, 1.511694908494142)

$rows = 2..25
foreach ($col in $colData.Keys) {
    $values = $colData[$col]
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $rows[$i]
        $ws.Range("$col$r").Value = $values[$i]
    }
}

Write-Output "Updated $($colData.Keys.Count) columns across $($rows.Count) rows."
